# "Wireframes version 2." -> "Wireframes version 1."
# (document text is "Version 2." -> "Version 1.")
#
# The paragraph text is built from several runs:
#   "Versi" | "on" | proofErr(spellEnd) | " 2" | bookmark | "."
# and needs to become:
#   "Version" | proofErr(spellEnd) | " 1." | bookmark
#
# i.e. the two "Version" runs merge into one, the " 2" run becomes " 1.",
# and the trailing "." run disappears (its "." moves into the " 1." run),
# while the _GoBack bookmark in between is left untouched.

$d = $word.ActiveDocument

# Step 1: turn the "2" in the " 2" run into "1." (in place, single run)
$r1 = $d.Range(8, 9)
$r1.Text = "1."

# Step 2: remove the now-redundant trailing "." run
$r2 = $d.Range(10, 11)
$r2.Text = ""

# Step 3: merge the "Versi" + "on" runs into a single "Version" run.
# (Setting identical text is a no-op for Word, so nudge it through a
# differing placeholder first, then fix the casing back.)
$r3 = $d.Range(0, 7)
$r3.Text = "VERSION"
$r3 = $d.Range(0, 7)
$r3.Text = "Version"
